$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price (D) and Volume (E) columns for the rows
# we touch, so values such as "314.90" or "0.846" are not silently
# reinterpreted as numbers (and lose trailing zeros) by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("D2").Value = '42.363.19'
$ws.Range("E2").Value = '  -1.29%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("D3").Value = '2.525.02'
$ws.Range("E3").Value = '  -0.70%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("D5").Value = '314.90'
$ws.Range("E5").Value = '  +3.08%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("D6").Value = '93.98'
$ws.Range("E6").Value = '  -5.40%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("D7").Value = '0.573'
$ws.Range("E7").Value = '  -0.91%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.03%  '

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -3.30%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("D10").Value = '35.41'
$ws.Range("E10").Value = '  -4.79%  '

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -2.31%  '

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -1.39%  '

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.41%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("D14").Value = '2.910.53'
$ws.Range("E14").Value = '  -0.66%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("D15").Value = '15.43'
$ws.Range("E15").Value = '  +0.96%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("D16").Value = '2.506.03'
$ws.Range("E16").Value = '  -2.74%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("D17").Value = '0.846'
$ws.Range("E17").Value = '  -3.42%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("D18").Value = '42.422.87'
$ws.Range("E18").Value = '  -1.18%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("D19").Value = '12.83'
$ws.Range("E19").Value = '  -2.53%  '

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.01%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0957'

$ws.Range("D22").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("D22").Value = '70.53'
$ws.Range("E22").Value = '  -1.64%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("D23").Value = '249.63'
$ws.Range("E23").Value = '  -1.91%  '

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.61%  '

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -3.03%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("D26").Value = '26.57'
$ws.Range("E26").Value = '  -4.25%  '

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.17%  '

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +3.73%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("B29").Value = 'InjectiveProtocol'
$ws.Range("C29").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D29").Value = '38.80'
$ws.Range("E29").Value = '  -0.03%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("B30").Value = 'Cosmos'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D30").Value = '10.08'
$ws.Range("E30").Value = '  -1.29%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("D31").Value = '5.90'
$ws.Range("E31").Value = '  -4.31%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("D32").Value = '155.98'
$ws.Range("E32").Value = '  -1.07%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("D33").Value = '2.11'
$ws.Range("E33").Value = '  -1.26%  '

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.04%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("D35").Value = '18.93'
$ws.Range("E35").Value = '  -0.20%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0779'
$ws.Range("E36").Value = '  -2.97%  '

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.92%  '

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -4.54%  '

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.55%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("D40").Value = '23.64'
$ws.Range("E40").Value = '  -2.85%  '

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.30%  '

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -2.96%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("D44").Value = '3.29'
$ws.Range("E44").Value = '  -5.22%  '

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -2.12%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("D46").Value = '2.014.52'
$ws.Range("E46").Value = '  -3.00%  '

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -2.57%  '

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -2.69%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("D49").Value = '2.764.14'
$ws.Range("E49").Value = '  -0.81%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("D50").Value = '101.75'
$ws.Range("E50").Value = '  -1.73%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("D51").Value = '72.28'
$ws.Range("E51").Value = '  -2.07%  '
